# Generate Report for Handoff
# Updates the "ab1afdcc-c517-45d0-be34-3eb793a40db6" row with its latest
# handoff/handback timestamps on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D7").Value = "2016-30-17 16:30:25"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E7").Value = "2016-03-17 16:30:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E7").Value = "2016-03-17 16:30:25"
